# add_column_in_sheet_differently_sorted: extend the "group" columns so the
# label choice reaches 7 more columns to the right (EH -> EO), mirroring the
# last existing column (EH) into the new columns EI:EO for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column EH is #138, the 7 new columns EI..EO are #139-#145.
$lastExistingCol = 138
$firstNewCol = 139
$lastNewCol = 145
$firstRow = 2
$lastRow = 15

# Step 1: copy the formatting of the last existing column onto the new
# columns in one shot, so the new cells are created (even the ones that stay
# blank) with the same style as the column they extend.
$src = $ws.Range($ws.Cells.Item($firstRow, $lastExistingCol), $ws.Cells.Item($lastRow, $lastExistingCol))
$src.Copy()
$dst = $ws.Range($ws.Cells.Item($firstRow, $firstNewCol), $ws.Cells.Item($lastRow, $lastNewCol))
$dst.PasteSpecial(-4122)

# Step 2: fill in the actual values, copying the last existing column's
# value across each of the new columns for every row (skip blank rows so we
# don't wipe out the blank cell created above).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $v = $ws.Cells.Item($r, $lastExistingCol).Formula
    if ($v -ne "") {
        for ($c = $firstNewCol; $c -le $lastNewCol; $c++) {
            $ws.Cells.Item($r, $c).Value = $v
        }
    }
}
